# 0. CTRL-E Autofill Example.xlsx — content update
#
# Replaces the old "CTRL-E Autofill Hack" title with updated flash-fill
# instructions, and appends guidance about creating a custom fill-handle
# list underneath the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the heading above the table (B2).
$ws.Range("B2").Value = "1. Use flash fill to complete these emails @xyz.com"

# Add the new instructional text below the table (rows 10, 11, 13).
$ws.Range("B10").Value = "2. Create a new custom list for the fill handle:"
$ws.Range("B11").Value = "Instructions are at https://support.microsoft.com/en-us/office/create-or-delete-a-custom-list-for-sorting-and-filling-data-d1cf624f-2d2b-44fa-814b-ba213ec2fd61"
$ws.Range("B13").Value = "Examples: Small/Medium/Large, North/East/South/West, etc."
